# "Separado o código de exportação da análise"
#
# Renames/reorders the three worksheets:
#   DadosCarga   -> Analise       (gains comparison columns F:P)
#   DadosSistema -> DadosCarga    (loses its "chaveCPF&Nome" column, rows reordered)
#   Analise      -> DadosSistema  (gains more "chave*" columns, rows reordered)
#
# and rewrites each sheet's data to match the new role.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rotate the three sheet names (3-cycle, via a temp name to avoid clashes)
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("DadosCarga").Name   = "TmpRenameSheet"
$wb.Worksheets.Item("DadosSistema").Name = "DadosCarga"
$wb.Worksheets.Item("Analise").Name      = "DadosSistema"
$wb.Worksheets.Item("TmpRenameSheet").Name = "Analise"

$wsAnalise      = $wb.Worksheets.Item("Analise")
$wsDadosCarga   = $wb.Worksheets.Item("DadosCarga")
$wsDadosSistema = $wb.Worksheets.Item("DadosSistema")

# ---------------------------------------------------------------------------
# 2) Analise sheet (ex-DadosCarga): columns A:E already correct (same row
#    order Mário / Anthony / Yasmin / Marina). Add comparison columns F:P.
# ---------------------------------------------------------------------------
$wsAnalise.Range("A1").Copy()
$wsAnalise.Range("F1:P1").PasteSpecial(-4122)

$wsAnalise.Cells.Item(1,6).Value  = "ComparaCPF"
$wsAnalise.Cells.Item(1,7).Value  = "chaveCPF&Nome"
$wsAnalise.Cells.Item(1,8).Value  = "ComparaCPF_Nome"
$wsAnalise.Cells.Item(1,9).Value  = "chaveCPF&DtNasc"
$wsAnalise.Cells.Item(1,10).Value = "ComparaCPF_DtNasc"
$wsAnalise.Cells.Item(1,11).Value = "chaveCPF&Sexo"
$wsAnalise.Cells.Item(1,12).Value = "ComparaCPF_Sexo"
$wsAnalise.Cells.Item(1,13).Value = "chaveCPF&Matricula"
$wsAnalise.Cells.Item(1,14).Value = "ComparaCPF_Matricula"
$wsAnalise.Cells.Item(1,15).Value = "chave"
$wsAnalise.Cells.Item(1,16).Value = "Resultado"

$wsAnalise.Cells.Item(2,6).Value  = "Localizado"
$wsAnalise.Cells.Item(2,7).Value  = "039.173.221-88Mário Manoel Calebe Moura"
$wsAnalise.Cells.Item(2,8).Value  = "Localizado"
$wsAnalise.Cells.Item(2,9).Value  = "039.173.221-881954-09-26"
$wsAnalise.Cells.Item(2,10).Value = "Localizado"
$wsAnalise.Cells.Item(2,11).Value = "039.173.221-88M"
$wsAnalise.Cells.Item(2,12).Value = "Localizado"
$wsAnalise.Cells.Item(2,13).Value = "039.173.221-8870297"
$wsAnalise.Cells.Item(2,14).Value = "Localizado"
$wsAnalise.Cells.Item(2,15).Value = "039.173.221-88Mário Manoel Calebe Moura1954-09-26M70297"
$wsAnalise.Cells.Item(2,16).Value = "Localizado"

$wsAnalise.Cells.Item(3,6).Value  = "Localizado"
$wsAnalise.Cells.Item(3,7).Value  = "876.995.634-09Anthony Henrique Costa"
$wsAnalise.Cells.Item(3,8).Value  = "Localizado"
$wsAnalise.Cells.Item(3,9).Value  = "876.995.634-091992-03-26"
$wsAnalise.Cells.Item(3,10).Value = "Localizado"
$wsAnalise.Cells.Item(3,11).Value = "876.995.634-09M"
$wsAnalise.Cells.Item(3,12).Value = "Localizado"
$wsAnalise.Cells.Item(3,13).Value = "876.995.634-0969314"
$wsAnalise.Cells.Item(3,14).Value = "Localizado"
$wsAnalise.Cells.Item(3,15).Value = "876.995.634-09Anthony Henrique Costa1992-03-26M69314"
$wsAnalise.Cells.Item(3,16).Value = "Localizado"

$wsAnalise.Cells.Item(4,6).Value  = "Localizado"
$wsAnalise.Cells.Item(4,7).Value  = "193.703.911-00Yasmin Eliane Agatha Gonçalves"
$wsAnalise.Cells.Item(4,8).Value  = "Localizado"
$wsAnalise.Cells.Item(4,9).Value  = "193.703.911-001945-07-17"
$wsAnalise.Cells.Item(4,10).Value = "Localizado"
$wsAnalise.Cells.Item(4,11).Value = "193.703.911-00F"
$wsAnalise.Cells.Item(4,12).Value = "Localizado"
$wsAnalise.Cells.Item(4,13).Value = "193.703.911-0061902"
$wsAnalise.Cells.Item(4,14).Value = "Localizado"
$wsAnalise.Cells.Item(4,15).Value = "193.703.911-00Yasmin Eliane Agatha Gonçalves1945-07-17F61902"
$wsAnalise.Cells.Item(4,16).Value = "Localizado"

$wsAnalise.Cells.Item(5,6).Value  = "Localizado"
$wsAnalise.Cells.Item(5,7).Value  = "369.960.476-41Marina Sophie Marlene da Luz"
$wsAnalise.Cells.Item(5,8).Value  = "Localizado"
$wsAnalise.Cells.Item(5,9).Value  = "369.960.476-411944-02-05"
$wsAnalise.Cells.Item(5,10).Value = "Localizado"
$wsAnalise.Cells.Item(5,11).Value = "369.960.476-41F"
$wsAnalise.Cells.Item(5,12).Value = "Localizado"
$wsAnalise.Cells.Item(5,13).Value = "369.960.476-4138045"
$wsAnalise.Cells.Item(5,14).Value = "Localizado"
$wsAnalise.Cells.Item(5,15).Value = "369.960.476-41Marina Sophie Marlene da Luz1944-02-05F38045"
$wsAnalise.Cells.Item(5,16).Value = "Localizado"

# ---------------------------------------------------------------------------
# 3) DadosCarga sheet (ex-DadosSistema): drop the "chaveCPF&Nome" column and
#    put the two middle rows (Mário / Anthony) back in load order.
# ---------------------------------------------------------------------------
$wsDadosCarga.Range("F1:F5").Clear()

$wsDadosCarga.Cells.Item(2,1).Value = "Mário Manoel Calebe Moura"
$wsDadosCarga.Cells.Item(2,2).Value = "039.173.221-88"
$wsDadosCarga.Cells.Item(2,3).Value = 19993
$wsDadosCarga.Cells.Item(2,4).Value = "M"
$wsDadosCarga.Cells.Item(2,5).Value = 70297

$wsDadosCarga.Cells.Item(3,1).Value = "Anthony Henrique Costa"
$wsDadosCarga.Cells.Item(3,2).Value = "876.995.634-09"
$wsDadosCarga.Cells.Item(3,3).Value = 33689
$wsDadosCarga.Cells.Item(3,4).Value = "M"
$wsDadosCarga.Cells.Item(3,5).Value = 69314

# ---------------------------------------------------------------------------
# 4) DadosSistema sheet (ex-Analise): reorder the Mário/Anthony rows, replace
#    the "ComparaCPF*" columns with extra "chave*" columns.
# ---------------------------------------------------------------------------
$wsDadosSistema.Range("A1").Copy()
$wsDadosSistema.Range("I1:J1").PasteSpecial(-4122)

$wsDadosSistema.Cells.Item(1,6).Value  = "chaveCPF&Nome"
$wsDadosSistema.Cells.Item(1,7).Value  = "chaveCPF&DtNasc"
$wsDadosSistema.Cells.Item(1,8).Value  = "chaveCPF&Sexo"
$wsDadosSistema.Cells.Item(1,9).Value  = "chaveCPF&Matricula"
$wsDadosSistema.Cells.Item(1,10).Value = "chave"

$wsDadosSistema.Cells.Item(2,1).Value = "Anthony Henrique Costa"
$wsDadosSistema.Cells.Item(2,2).Value = "876.995.634-09"
$wsDadosSistema.Cells.Item(2,3).Value = 33689
$wsDadosSistema.Cells.Item(2,4).Value = "M"
$wsDadosSistema.Cells.Item(2,5).Value = 69314
$wsDadosSistema.Cells.Item(2,6).Value = "876.995.634-09Anthony Henrique Costa"
$wsDadosSistema.Cells.Item(2,7).Value = "876.995.634-091992-03-26"
$wsDadosSistema.Cells.Item(2,8).Value = "876.995.634-09M"
$wsDadosSistema.Cells.Item(2,9).Value = "876.995.634-0969314"
$wsDadosSistema.Cells.Item(2,10).Value = "876.995.634-09Anthony Henrique Costa1992-03-26M69314"

$wsDadosSistema.Cells.Item(3,1).Value = "Mário Manoel Calebe Moura"
$wsDadosSistema.Cells.Item(3,2).Value = "039.173.221-88"
$wsDadosSistema.Cells.Item(3,3).Value = 19993
$wsDadosSistema.Cells.Item(3,4).Value = "M"
$wsDadosSistema.Cells.Item(3,5).Value = 70297
$wsDadosSistema.Cells.Item(3,6).Value = "039.173.221-88Mário Manoel Calebe Moura"
$wsDadosSistema.Cells.Item(3,7).Value = "039.173.221-881954-09-26"
$wsDadosSistema.Cells.Item(3,8).Value = "039.173.221-88M"
$wsDadosSistema.Cells.Item(3,9).Value = "039.173.221-8870297"
$wsDadosSistema.Cells.Item(3,10).Value = "039.173.221-88Mário Manoel Calebe Moura1954-09-26M70297"

$wsDadosSistema.Cells.Item(4,6).Value = "193.703.911-00Yasmin Eliane Agatha Gonçalves"
$wsDadosSistema.Cells.Item(4,7).Value = "193.703.911-001945-07-17"
$wsDadosSistema.Cells.Item(4,8).Value = "193.703.911-00F"
$wsDadosSistema.Cells.Item(4,9).Value = "193.703.911-0061902"
$wsDadosSistema.Cells.Item(4,10).Value = "193.703.911-00Yasmin Eliane Agatha Gonçalves1945-07-17F61902"

$wsDadosSistema.Cells.Item(5,6).Value = "369.960.476-41Marina Sophie Marlene da Luz"
$wsDadosSistema.Cells.Item(5,7).Value = "369.960.476-411944-02-05"
$wsDadosSistema.Cells.Item(5,8).Value = "369.960.476-41F"
$wsDadosSistema.Cells.Item(5,9).Value = "369.960.476-4138045"
$wsDadosSistema.Cells.Item(5,10).Value = "369.960.476-41Marina Sophie Marlene da Luz1944-02-05F38045"
